# Fruta / hortaliza, semanal
# Insert 4 new weekly records at the top of the "Palta" price block (rows 201-204),
# pushing the existing rows 201-293 down to 205-297.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 201 (existing rows 201:293 shift to 205:297).
$ws.Rows("201:204").Insert()

# Common fields shared by every record in this block.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"

# New row 201: Fuerte / Especial
$r = 201
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 45205
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fuerte"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 208
$ws.Cells.Item($r, 14).Value = 35000
$ws.Cells.Item($r, 15).Value = 36000
$ws.Cells.Item($r, 16).Value = 35500
$ws.Cells.Item($r, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Perú"
$ws.Cells.Item($r, 19).Value = 3550
$ws.Cells.Item($r, 20).Value = 10

# New row 202: Fuerte / Primera
$r = 202
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 45205
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fuerte"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 208
$ws.Cells.Item($r, 14).Value = 33000
$ws.Cells.Item($r, 15).Value = 34000
$ws.Cells.Item($r, 16).Value = 33500
$ws.Cells.Item($r, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Perú"
$ws.Cells.Item($r, 19).Value = 3350
$ws.Cells.Item($r, 20).Value = 10

# New row 203: Fuerte / Segunda
$r = 203
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 45205
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fuerte"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 208
$ws.Cells.Item($r, 14).Value = 31000
$ws.Cells.Item($r, 15).Value = 32000
$ws.Cells.Item($r, 16).Value = 31500
$ws.Cells.Item($r, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Perú"
$ws.Cells.Item($r, 19).Value = 3150
$ws.Cells.Item($r, 20).Value = 10

# New row 204: Fuerte / Tercera
$r = 204
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 45205
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fuerte"
$ws.Cells.Item($r, 12).Value = "Tercera"
$ws.Cells.Item($r, 13).Value = 208
$ws.Cells.Item($r, 14).Value = 29000
$ws.Cells.Item($r, 15).Value = 30000
$ws.Cells.Item($r, 16).Value = 29500
$ws.Cells.Item($r, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Perú"
$ws.Cells.Item($r, 19).Value = 2950
$ws.Cells.Item($r, 20).Value = 10

# Make sure the D column (date) keeps its date style/format on the new rows.
$ws.Range("D201:D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
